# Update cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must stay TEXT (matches original inlineStr
# cells) are written via NumberFormat "@" so Excel does not coerce them to numbers,
# then the format is reset back to Normal so no stray style is left behind.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '67.471.89'
$ws.Range("E2").Value = '  +4.26%  '
$ws.Range("D3").Value = '3.248.44'
$ws.Range("E3").Value = '  +2.67%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.11%  '
Set-TextValue $ws.Range("D5") '578.24'
$ws.Range("E5").Value = '  +2.34%  '
Set-TextValue $ws.Range("D6") '181.57'
$ws.Range("E6").Value = '  +6.19%  '
Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -4.68%  '
$ws.Range("D9").Value = '3.243.03'
$ws.Range("E9").Value = '  +2.53%  '
Set-TextValue $ws.Range("D10") '0.130'
$ws.Range("E10").Value = '  +4.30%  '
Set-TextValue $ws.Range("D11") '6.79'
$ws.Range("E11").Value = '  +3.30%  '
Set-TextValue $ws.Range("D12") '0.414'
$ws.Range("E12").Value = '  +5.09%  '
$ws.Range("D13").Value = '3.796.25'
$ws.Range("E13").Value = '  +1.91%  '
Set-TextValue $ws.Range("D14") '0.138'
$ws.Range("E14").Value = '  +1.61%  '
Set-TextValue $ws.Range("D15") '28.03'
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("D16").Value = '67.411.20'
$ws.Range("E16").Value = '  +4.29%  '
Set-TextValue $ws.Range("D17") '0.0000168'
$ws.Range("E17").Value = '  +2.78%  '
$ws.Range("D18").Value = '3.239.06'
$ws.Range("E18").Value = '  +2.11%  '
Set-TextValue $ws.Range("D19") '5.81'
$ws.Range("E19").Value = '  +1.72%  '
Set-TextValue $ws.Range("D20") '13.44'
$ws.Range("E20").Value = '  +3.64%  '
Set-TextValue $ws.Range("D21") '375.08'
$ws.Range("E21").Value = '  +5.55%  '
Set-TextValue $ws.Range("D22") '7.59'
$ws.Range("E22").Value = '  +4.91%  '
$ws.Range("E23").Value = '  -0.32%  '
Set-TextValue $ws.Range("D24") '71.17'
$ws.Range("E24").Value = '  +3.70%  '
Set-TextValue $ws.Range("D25") '0.510'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("E26").Value = '  +1.66%  '
Set-TextValue $ws.Range("D27") '9.63'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  +3.14%  '
$ws.Range("E29").Value = '  +0.34%  '
Set-TextValue $ws.Range("D30") '1.97'
$ws.Range("E30").Value = '  +3.83%  '
Set-TextValue $ws.Range("D31") '5.69'
$ws.Range("E31").Value = '  +5.56%  '
Set-TextValue $ws.Range("D32") '22.62'
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("E34").Value = '  +4.64%  '
Set-TextValue $ws.Range("D35") '6.86'
$ws.Range("E35").Value = '  +2.79%  '
Set-TextValue $ws.Range("D36") '162.14'
$ws.Range("E36").Value = '  +5.76%  '
Set-TextValue $ws.Range("D37") '1.49'
$ws.Range("E37").Value = '  +3.22%  '
Set-TextValue $ws.Range("D38") '0.857'
$ws.Range("E38").Value = '  +3.01%  '
Set-TextValue $ws.Range("D39") '1.86'
$ws.Range("E39").Value = '  +7.15%  '
Set-TextValue $ws.Range("D40") '6.79'
$ws.Range("E40").Value = '  +12.38%  '
Set-TextValue $ws.Range("D41") '26.64'
$ws.Range("E41").Value = '  +1.55%  '
Set-TextValue $ws.Range("D42") '2.61'
$ws.Range("E42").Value = '  +3.13%  '
Set-TextValue $ws.Range("D43") '363.37'
$ws.Range("E43").Value = '  +12.37%  '
Set-TextValue $ws.Range("D44") '4.45'
$ws.Range("E44").Value = '  +6.19%  '
$ws.Range("D45").Value = '2.718.16'
$ws.Range("E45").Value = '  +2.43%  '
Set-TextValue $ws.Range("D46") '25.60'
$ws.Range("E46").Value = '  +5.61%  '
Set-TextValue $ws.Range("D47") '40.57'
$ws.Range("E47").Value = '  +3.57%  '
Set-TextValue $ws.Range("D48") '0.0672'
$ws.Range("E48").Value = '  +2.94%  '
Set-TextValue $ws.Range("D49") '0.0279'
$ws.Range("E49").Value = '  +2.19%  '
Set-TextValue $ws.Range("D50") '0.995'
$ws.Range("E50").Value = '  +6.10%  '
$ws.Range("E51").Value = '  -0.22%  '
